# Auto-generated edit script: update Leve profit/price columns (H-N) for rows
# changed by the scheduled market-data refresh, across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 290.58334
$ws.Range("I33").Value = 199.5
$ws.Range("K33").Value = 199.5
$ws.Range("M33").Value = 29.5

$ws.Range("H53").Value = 270
$ws.Range("I53").Value = 274
$ws.Range("K53").Value = 274
$ws.Range("M53").Value = 363

$ws.Range("H64").Value = 2744.2778
$ws.Range("I64").Value = 2714
$ws.Range("J64").Value = 2763.5454
$ws.Range("K64").Value = 2714
$ws.Range("L64").Value = 2763.5454
$ws.Range("M64").Value = -2466
$ws.Range("N64").Value = -3259.5454

$ws.Range("H67").Value = 2744.2778
$ws.Range("I67").Value = 2714
$ws.Range("J67").Value = 2763.5454
$ws.Range("K67").Value = 2714
$ws.Range("L67").Value = 2763.5454
$ws.Range("M67").Value = -1856
$ws.Range("N67").Value = -4479.5454

$ws.Range("H129").Value = 918
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 1012.1818
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 3036.5454
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -13036.5454

$ws.Range("H132").Value = 19309628
$ws.Range("I132").Value = 20491308
$ws.Range("K132").Value = 61473924
$ws.Range("M132").Value = -61471394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 780.2222
$ws.Range("I2").Value = 797.82355
$ws.Range("J2").Value = 750.3
$ws.Range("K2").Value = 797.82355
$ws.Range("L2").Value = 750.3
$ws.Range("M2").Value = -684.82355
$ws.Range("N2").Value = -976.3

$ws.Range("H34").Value = 14933.333
$ws.Range("I34").Value = 3000
$ws.Range("J34").Value = 38800
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 38800
$ws.Range("M34").Value = -2729
$ws.Range("N34").Value = -39342

$ws.Range("H116").Value = 780.2222
$ws.Range("I116").Value = 797.82355
$ws.Range("J116").Value = 750.3
$ws.Range("K116").Value = 797.82355
$ws.Range("L116").Value = 750.3
$ws.Range("M116").Value = 1496.17645
$ws.Range("N116").Value = -5338.3

$ws.Range("H122").Value = 2231.7856
$ws.Range("I122").Value = 1249.6154
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 3748.8462
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -1298.8462
$ws.Range("N122").Value = -49900

$ws.Range("H132").Value = 2353
$ws.Range("I132").Value = 1730.1132
$ws.Range("K132").Value = 5190.3396
$ws.Range("M132").Value = -2660.3396

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 780.2222
$ws.Range("I3").Value = 797.82355
$ws.Range("J3").Value = 750.3
$ws.Range("K3").Value = 797.82355
$ws.Range("L3").Value = 750.3
$ws.Range("M3").Value = -683.82355
$ws.Range("N3").Value = -978.3

$ws.Range("H22").Value = 644.6111
$ws.Range("I22").Value = 425.1875
$ws.Range("K22").Value = 425.1875
$ws.Range("M22").Value = -252.1875

$ws.Range("H64").Value = 344.8889
$ws.Range("I64").Value = 235.33333
$ws.Range("J64").Value = 399.66666
$ws.Range("K64").Value = 235.33333
$ws.Range("L64").Value = 399.66666
$ws.Range("M64").Value = -10.33332999999999
$ws.Range("N64").Value = -849.66666

$ws.Range("H67").Value = 344.8889
$ws.Range("I67").Value = 235.33333
$ws.Range("J67").Value = 399.66666
$ws.Range("K67").Value = 235.33333
$ws.Range("L67").Value = 399.66666
$ws.Range("M67").Value = 544.6666700000001
$ws.Range("N67").Value = -1959.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2245.0164
$ws.Range("I31").Value = 912.8857400000001
$ws.Range("J31").Value = 4038.2693
$ws.Range("K31").Value = 912.8857400000001
$ws.Range("L31").Value = 4038.2693
$ws.Range("M31").Value = -617.8857400000001
$ws.Range("N31").Value = -4628.2693

$ws.Range("H34").Value = 2245.0164
$ws.Range("I34").Value = 912.8857400000001
$ws.Range("J34").Value = 4038.2693
$ws.Range("K34").Value = 912.8857400000001
$ws.Range("L34").Value = 4038.2693
$ws.Range("M34").Value = -710.8857400000001
$ws.Range("N34").Value = -4442.2693

$ws.Range("H58").Value = 1924.5922
$ws.Range("I58").Value = 1601.4849
$ws.Range("J58").Value = 4057.1
$ws.Range("K58").Value = 1601.4849
$ws.Range("L58").Value = 4057.1
$ws.Range("M58").Value = -1398.4849
$ws.Range("N58").Value = -4463.1

$ws.Range("H94").Value = 2115.1333
$ws.Range("I94").Value = 1747.5
$ws.Range("J94").Value = 2171.6924
$ws.Range("K94").Value = 1747.5
$ws.Range("L94").Value = 2171.6924
$ws.Range("M94").Value = -1296.5
$ws.Range("N94").Value = -3073.6924

$ws.Range("H99").Value = 6727.9375
$ws.Range("I99").Value = 5302.1113
$ws.Range("J99").Value = 8561.143
$ws.Range("K99").Value = 5302.1113
$ws.Range("L99").Value = 8561.143
$ws.Range("M99").Value = -3804.1113
$ws.Range("N99").Value = -11557.143

$ws.Range("H126").Value = 6727.9375
$ws.Range("I126").Value = 5302.1113
$ws.Range("J126").Value = 8561.143
$ws.Range("K126").Value = 15906.3339
$ws.Range("L126").Value = 25683.429
$ws.Range("M126").Value = -13436.3339
$ws.Range("N126").Value = -30623.429

$ws.Range("H132").Value = 2977.9788
$ws.Range("I132").Value = 1952.5588
$ws.Range("J132").Value = 5659.846
$ws.Range("K132").Value = 5857.6764
$ws.Range("L132").Value = 16979.538
$ws.Range("M132").Value = -3327.6764
$ws.Range("N132").Value = -22039.538

$ws.Range("H134").Value = 4436.75
$ws.Range("I134").Value = 4863.52
$ws.Range("J134").Value = 3466.818
$ws.Range("K134").Value = 14590.56
$ws.Range("L134").Value = 10400.454
$ws.Range("M134").Value = -12055.56
$ws.Range("N134").Value = -15470.454

$ws.Range("H136").Value = 1924.5922
$ws.Range("I136").Value = 1601.4849
$ws.Range("J136").Value = 4057.1
$ws.Range("K136").Value = 4804.4547
$ws.Range("L136").Value = 12171.3
$ws.Range("M136").Value = -2254.4547
$ws.Range("N136").Value = -17271.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1109.9048
$ws.Range("I5").Value = 374.82858
$ws.Range("J5").Value = 4785.2856
$ws.Range("K5").Value = 1124.48574
$ws.Range("L5").Value = 14355.8568
$ws.Range("M5").Value = -1012.48574
$ws.Range("N5").Value = -14579.8568

$ws.Range("H113").Value = 585.58185
$ws.Range("I113").Value = 586.1539
$ws.Range("J113").Value = 584.1875
$ws.Range("K113").Value = 1758.4617
$ws.Range("L113").Value = 1752.5625
$ws.Range("M113").Value = 411.5382999999999
$ws.Range("N113").Value = -6092.5625

$ws.Range("H132").Value = 3181.1724
$ws.Range("I132").Value = 820.8889
$ws.Range("K132").Value = 7388.0001
$ws.Range("M132").Value = -4858.0001

$ws.Range("H135").Value = 1109.9048
$ws.Range("I135").Value = 374.82858
$ws.Range("J135").Value = 4785.2856
$ws.Range("K135").Value = 3373.45722
$ws.Range("L135").Value = 43067.5704
$ws.Range("M135").Value = -838.4572199999998
$ws.Range("N135").Value = -48137.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3790.798
$ws.Range("I126").Value = 2714.8167
$ws.Range("K126").Value = 8144.4501
$ws.Range("M126").Value = -5674.4501

$ws.Range("H132").Value = 2441.6
$ws.Range("I132").Value = 1731.1305
$ws.Range("J132").Value = 3803.3333
$ws.Range("K132").Value = 5193.3915
$ws.Range("L132").Value = 11409.9999
$ws.Range("M132").Value = -2663.3915
$ws.Range("N132").Value = -16469.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 14000
$ws.Range("I29").Value = 11000
$ws.Range("J29").Value = 20000
$ws.Range("K29").Value = 11000
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = -10705
$ws.Range("N29").Value = -20590

$ws.Range("H40").Value = 5863.1875
$ws.Range("J40").Value = 5301
$ws.Range("L40").Value = 5301
$ws.Range("N40").Value = -5573

$ws.Range("H61").Value = 1440.1428
$ws.Range("I61").Value = 1437.0625
$ws.Range("J61").Value = 1450
$ws.Range("K61").Value = 1437.0625
$ws.Range("L61").Value = 1450
$ws.Range("M61").Value = -1235.0625
$ws.Range("N61").Value = -1854

$ws.Range("H100").Value = 2034.8
$ws.Range("I100").Value = 1847.2222
$ws.Range("J100").Value = 2316.1667
$ws.Range("K100").Value = 1847.2222
$ws.Range("L100").Value = 2316.1667
$ws.Range("M100").Value = -1306.2222
$ws.Range("N100").Value = -3398.1667

$ws.Range("H113").Value = 1440.1428
$ws.Range("I113").Value = 1437.0625
$ws.Range("J113").Value = 1450
$ws.Range("K113").Value = 1437.0625
$ws.Range("L113").Value = 1450
$ws.Range("M113").Value = 732.9375
$ws.Range("N113").Value = -5790

$ws.Range("H122").Value = 4289.7744
$ws.Range("I122").Value = 3295.36
$ws.Range("J122").Value = 8433.166999999999
$ws.Range("K122").Value = 9886.08
$ws.Range("L122").Value = 25299.501
$ws.Range("M122").Value = -7436.08
$ws.Range("N122").Value = -30199.501

$ws.Range("H132").Value = 4918.173
$ws.Range("I132").Value = 1633.7826
$ws.Range("K132").Value = 4901.3478
$ws.Range("M132").Value = -2371.3478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 35715420
$ws.Range("I81").Value = 35715420
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 71430840
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -71429779
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 35715420
$ws.Range("I84").Value = 35715420
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 357154200
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -357148896
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value = 264.5
$ws.Range("I113").Value = 245.94118
$ws.Range("J113").Value = 299.55554
$ws.Range("K113").Value = 737.82354
$ws.Range("L113").Value = 898.66662
$ws.Range("M113").Value = 1432.17646
$ws.Range("N113").Value = -5238.66662

$ws.Range("H132").Value = 5377813.5
$ws.Range("I132").Value = 514.3570999999999
$ws.Range("K132").Value = 1543.0713
$ws.Range("M132").Value = 986.9287000000002
